$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the neighbouring header cell (AC1) onto the new
# header cells so they match the rest of the header row (bold, bordered,
# centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record columns for every data row (2 through 48)
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 56   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 106  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
